# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Updates the cryptos price-tracker sheet with refreshed values (Coin, Link,
# Price, Volume(1h)) as captured by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '56.665.27'
$ws.Range("D3").Value2 = '2.985.16'
$ws.Range("E3").Value2 = '  -5.39%  '
$ws.Range("E4").Value2 = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '497.02'
$ws.Range("E5").Value2 = '  -5.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '135.47'
$ws.Range("E6").Value2 = '  +0.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.425'
$ws.Range("E8").Value2 = '  -4.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '7.22'
$ws.Range("E9").Value2 = '  -1.16%  '
$ws.Range("E10").Value2 = '  -3.35%  '
$ws.Range("E11").Value2 = '  -7.15%  '
$ws.Range("E12").Value2 = '  -0.73%  '
$ws.Range("D13").Value2 = '3.490.31'
$ws.Range("E13").Value2 = '  -5.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '25.29'
$ws.Range("E14").Value2 = '  -0.55%  '
$ws.Range("D15").Value2 = '56.581.60'
$ws.Range("E15").Value2 = '  -3.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '0.0000148'
$ws.Range("E16").Value2 = '  -2.81%  '
$ws.Range("D17").Value2 = '2.980.40'
$ws.Range("E17").Value2 = '  -5.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '5.80'
$ws.Range("E18").Value2 = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '12.36'
$ws.Range("E19").Value2 = '  -5.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '7.80'
$ws.Range("E20").Value2 = '  -2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '325.60'
$ws.Range("E21").Value2 = '  -5.03%  '
$ws.Range("E22").Value2 = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '0.470'
$ws.Range("E23").Value2 = '  -7.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '61.74'
$ws.Range("E24").Value2 = '  -7.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '1.00'
$ws.Range("E25").Value2 = '  -0.01%  '
$ws.Range("E26").Value2 = '  -5.48%  '
$ws.Range("D27").Value2 = '0.0₃0897'
$ws.Range("E27").Value2 = '  -5.65%  '
$ws.Range("E28").Value2 = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '6.49'
$ws.Range("E29").Value2 = '  -5.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '6.75'
$ws.Range("E30").Value2 = '  -2.43%  '
$ws.Range("E31").Value2 = '  -6.95%  '
$ws.Range("B32").Value2 = 'EthereumClassic'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '20.32'
$ws.Range("E32").Value2 = '  -5.37%  '
$ws.Range("B33").Value2 = 'Fetch.AI'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '1.16'
$ws.Range("E33").Value2 = '  -7.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '152.62'
$ws.Range("E34").Value2 = '  -4.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '4.46'
$ws.Range("E35").Value2 = '  -8.03%  '
$ws.Range("E36").Value2 = '  -7.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '5.61'
$ws.Range("E37").Value2 = '  -10.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.0672'
$ws.Range("E38").Value2 = '  -2.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '23.28'
$ws.Range("E39").Value2 = '  -2.67%  '
$ws.Range("D40").Value2 = '3.011.32'
$ws.Range("E40").Value2 = '  -5.49%  '
$ws.Range("E41").Value2 = '  -9.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '1.00'
$ws.Range("E42").Value2 = '  +0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.640'
$ws.Range("E43").Value2 = '  -8.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.999'
$ws.Range("E44").Value2 = '  -8.26%  '
$ws.Range("B45").Value2 = 'Maker'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value2 = '2.228.64'
$ws.Range("E45").Value2 = '  -2.71%  '
$ws.Range("B46").Value2 = 'Stacks'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '1.41'
$ws.Range("E46").Value2 = '  -3.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '3.56'
$ws.Range("E47").Value2 = '  -9.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '1.95'
$ws.Range("E48").Value2 = '  +5.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.0236'
$ws.Range("E49").Value2 = '  +0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '5.74'
$ws.Range("E50").Value2 = '  -7.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '19.05'
$ws.Range("E51").Value2 = '  -8.13%  '
